$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '29.105.26'
Set-TextValue "E2" '  +0.12%  '

Set-TextValue "D3" '1.831.68'
Set-TextValue "E3" '  -0.22%  '

Set-TextValue "D4" '0.9988'
Set-TextValue "E4" '  -0.15%  '

Set-TextValue "D5" '242.95'
Set-TextValue "E5" '  +0.00%  '

Set-TextValue "D6" '0.6281'
Set-TextValue "E6" '  +0.22%  '

Set-TextValue "E7" '  -0.05%  '

Set-TextValue "D8" '0.07520'
Set-TextValue "E8" '  -0.98%  '

Set-TextValue "E9" '  -0.07%  '

Set-TextValue "D10" '23.22'
Set-TextValue "E10" '  +2.65%  '

Set-TextValue "D11" '0.07685'
Set-TextValue "E11" '  -0.77%  '

Set-TextValue "D12" '1.832.70'
Set-TextValue "E12" '  -0.59%  '

Set-TextValue "D13" '5.027'
Set-TextValue "E13" '  +1.17%  '

Set-TextValue "D14" '0.6692'
Set-TextValue "E14" '  +0.56%  '

Set-TextValue "D15" '82.87'
Set-TextValue "E15" '  -0.10%  '

Set-TextValue "D16" '0.000009381'
Set-TextValue "E16" '  -5.70%  '

Set-TextValue "D17" '5.992'
Set-TextValue "E17" '  -1.29%  '

Set-TextValue "D18" '29.108.44'
Set-TextValue "E18" '  +0.06%  '

Set-TextValue "D19" '2.078.59'
Set-TextValue "E19" '  -0.34%  '

Set-TextValue "E20" '  +1.47%  '

Set-TextValue "D21" '223.21'
Set-TextValue "E21" '  -1.76%  '

Set-TextValue "E22" '  +0.06%  '

Set-TextValue "D23" '7.145'
Set-TextValue "E23" '  -1.00%  '

Set-TextValue "D24" '1.000'
Set-TextValue "E24" '  -0.10%  '

Set-TextValue "D25" '160.09'
Set-TextValue "E25" '  +0.38%  '

Set-TextValue "D26" '0.1398'
Set-TextValue "E26" '  +1.02%  '

Set-TextValue "D27" '8.502'
Set-TextValue "E27" '  -0.18%  '

Set-TextValue "E28" '  -0.37%  '

Set-TextValue "D29" '1.495'
Set-TextValue "E29" '  -0.19%  '

Set-TextValue "D30" '0.05825'
Set-TextValue "E30" '  +10.87%  '

Set-TextValue "D31" '4.163'
Set-TextValue "E31" '  +1.31%  '

Set-TextValue "D32" '4.120'
Set-TextValue "E32" '  +2.44%  '

Set-TextValue "D33" '1.203'
Set-TextValue "E33" '  +0.80%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D34" '0.7414'
Set-TextValue "E34" '  +0.88%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D35" '1.834'
Set-TextValue "E35" '  -0.39%  '

Set-TextValue "E36" '  +0.09%  '

Set-TextValue "D37" '2.667'
Set-TextValue "E37" '  -0.89%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D38" '1.225.30'
Set-TextValue "E38" '  -1.15%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D39" '2.762'
Set-TextValue "E39" '  -0.05%  '

Set-TextValue "D40" '0.01781'
Set-TextValue "E40" '  -0.29%  '

Set-TextValue "D41" '6.496'
Set-TextValue "E41" '  +1.92%  '

Set-TextValue "D42" '0.8918'
Set-TextValue "E42" '  -0.85%  '

Set-TextValue "E43" '  -0.06%  '

Set-TextValue "D44" '101.93'
Set-TextValue "E44" '  -0.02%  '

Set-TextValue "E45" '  -0.46%  '

Set-TextValue "E46" '  +2.30%  '

Set-TextValue "E47" '  -2.30%  '

Set-TextValue "D48" '0.5087'
Set-TextValue "E48" '  -0.61%  '

Set-TextValue "D49" '0.07550'
Set-TextValue "E49" '  +12.26%  '

Set-TextValue "E50" '  +0.61%  '

Set-TextValue "D51" '8.995'
Set-TextValue "E51" '  +1.34%  '
